$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 3, columns T-W from 0.99 to 1
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 1
$ws.Range("V3").Value = 1
$ws.Range("W3").Value = 1

# Update the selection to match the new used range
$ws.Range("A1:X7").Select()
